$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Sheet "year" (3rd sheet) -- add the new scaling-rule rows (10-34).
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

# --- Row 10's note is typed first (creates the new shared string before any
#     of the new country codes so the shared-string table order matches the
#     authoring session). ---
$fsuNote  = "Do not scale FSU countries before 1992 since data are unlikely to match inventory assumptions"
$yugoNote = "Do not scale former Yugosloavia countries before 1992 since data are unlikely to match inventory assumptions"
$railNote = "Only scale rail transportation for a few years to avoid large jumps in EF"
$albNote  = "Avoid albania jump in EF after 2008 (CHECK: may need to also stop emep scaling)"

$ws3.Cells.Item(10, 8).Value = $fsuNote

# --- Former Soviet Union countries: rows 10-24 ---
$fsuCountries = @("arm","aze","blr","est","geo","kaz","kgz","ltu","lva","mda","rus","tjk","tkm","ukr","uzb")
for ($i = 0; $i -lt $fsuCountries.Length; $i++) {
    $ws3.Cells.Item(10 + $i, 1).Value = $fsuCountries[$i]
}

# --- Former Yugoslavia countries: rows 25-32 ---
$yugoCountries = @("scg","srb","svk","svn","bih","hrv","mkd","mne")
for ($i = 0; $i -lt $yugoCountries.Length; $i++) {
    $ws3.Cells.Item(25 + $i, 1).Value = $yugoCountries[$i]
}

# --- Fill in the rest of rows 10-24 (B-G), then H last for that block ---
for ($r = 10; $r -le 24; $r++) {
    $ws3.Cells.Item($r, 2).Value = "all"
    $ws3.Cells.Item($r, 3).Value = "NA"
    $ws3.Cells.Item($r, 4).Value = "NA"
    $ws3.Cells.Item($r, 5).Value = "NA"
    $ws3.Cells.Item($r, 6).Value = 1992
    $ws3.Cells.Item($r, 7).Value = 2020
    $ws3.Cells.Item($r, 8).Value = $fsuNote
}

# --- Fill in the rest of rows 25-32 (B-G) ---
for ($r = 25; $r -le 32; $r++) {
    $ws3.Cells.Item($r, 2).Value = "all"
    $ws3.Cells.Item($r, 3).Value = "NA"
    $ws3.Cells.Item($r, 4).Value = "NA"
    $ws3.Cells.Item($r, 5).Value = "NA"
    $ws3.Cells.Item($r, 6).Value = 1992
    $ws3.Cells.Item($r, 7).Value = 2020
}

# --- New Yugoslavia note, written after the A column fill-down, then spread
#     across H25:H32 (mirrors the authoring order recovered from the diff). ---
for ($r = 25; $r -le 32; $r++) {
    $ws3.Cells.Item($r, 8).Value = $yugoNote
}

# --- Row 33: "all" / rail-transportation sector rule ---
$ws3.Cells.Item(33, 1).Value = "all"
$ws3.Cells.Item(33, 2).Value = "1A3c"
$ws3.Cells.Item(33, 3).Value = "NA"
$ws3.Cells.Item(33, 4).Value = "NA"
$ws3.Cells.Item(33, 5).Value = "NA"
$ws3.Cells.Item(33, 6).Value = 2000
$ws3.Cells.Item(33, 7).Value = 2010
$ws3.Cells.Item(33, 8).Value = $railNote

# --- Row 34: Albania rule ---
$ws3.Cells.Item(34, 1).Value = "alb"
$ws3.Cells.Item(34, 2).Value = "1A3c"
$ws3.Cells.Item(34, 3).Value = "NA"
$ws3.Cells.Item(34, 4).Value = "NA"
$ws3.Cells.Item(34, 5).Value = "NA"
$ws3.Cells.Item(34, 6).Value = 2000
$ws3.Cells.Item(34, 7).Value = 2008
$ws3.Cells.Item(34, 8).Value = $albNote

# --- Style rows 33 & 34: column B and H get an explicit black font, and H
#     also gets the scientific-notation number format (the new cellXfs entry
#     that appears in the saved styles table). ---
$ws3.Range("B33").Font.Color = 0
$ws3.Range("B34").Font.Color = 0
$ws3.Range("H33").Font.Color = 0
$ws3.Range("H33").NumberFormat = "0.00E+00"
$ws3.Range("H34").Font.Color = 0
$ws3.Range("H34").NumberFormat = "0.00E+00"

# ---------------------------------------------------------------------------
# 2) View-state updates: sheet "map" loses focus/selection moves, sheet
#    "year" becomes the active tab with a frozen header row and a new
#    selection, and the workbook remembers "year" as the active sheet.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Activate()
$ws1.Range("B21:C21").Select()

$ws3.Activate()
$ws3.Range("A2").Select()
$win = $excel.ActiveWindow
$win.SplitRow = 1
$win.SplitColumn = 0
$win.FreezePanes = $true
$ws3.Range("G35").Select()
